$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F4").Value = 374
$ws.Range("F5").Value = 8292
$ws.Range("F7").Value = 112
$ws.Range("F8").Value = 2205
$ws.Range("F10").Value = 196
$ws.Range("F11").Value = 85
$ws.Range("F13").Value = 625
$ws.Range("F15").Value = 7246
$ws.Range("F16").Value = 450
$ws.Range("F18").Value = 7542
$ws.Range("F20").Value = 57189
$ws.Range("F21").Value = 4663
$ws.Range("F23").Value = 1051
$ws.Range("F24").Value = 919
$ws.Range("F25").Value = 477
$ws.Range("F27").Value = 906
$ws.Range("F29").Value = 610
$ws.Range("F30").Value = 4970
$ws.Range("F32").Value = 88
$ws.Range("F34").Value = 886
$ws.Range("F35").Value = 1293
$ws.Range("F36").Value = 1573
$ws.Range("F39").Value = 215
$ws.Range("F41").Value = 721
$ws.Range("F43").Value = 777
$ws.Range("F44").Value = 231
$ws.Range("F45").Value = 15
$ws.Range("F46").Value = 188
$ws.Range("F47").Value = 11
$ws.Range("F48").Value = 51

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F2").Value = 191
$ws.Range("F4").Value = 56
$ws.Range("F6").Value = 127
$ws.Range("F8").Value = 47
$ws.Range("F10").Value = 7561
$ws.Range("F20").Value = 20
$ws.Range("F22").Value = 2
$ws.Range("F40").Value = 109
$ws.Range("F41").Value = 185
$ws.Range("F47").Value = 269

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F5").Value = 1580
$ws.Range("F9").Value = 9400
$ws.Range("F15").Value = 244
$ws.Range("F16").Value = 2167
$ws.Range("F17").Value = 26
$ws.Range("F18").Value = 453

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F4").Value = 374
$ws.Range("F5").Value = 1580
$ws.Range("F11").Value = 112
$ws.Range("F12").Value = 244
$ws.Range("F13").Value = 2167
$ws.Range("F15").Value = 625
$ws.Range("F17").Value = 7246
$ws.Range("F18").Value = 57189
$ws.Range("F19").Value = 191
$ws.Range("F20").Value = 56
$ws.Range("F21").Value = 4663
$ws.Range("F22").Value = 1051
$ws.Range("F23").Value = 477
$ws.Range("F24").Value = 610
$ws.Range("F25").Value = 127
$ws.Range("F26").Value = 4970
$ws.Range("F28").Value = 88
$ws.Range("F29").Value = 886
$ws.Range("F30").Value = 1293
$ws.Range("F32").Value = 453
$ws.Range("F36").Value = 215
$ws.Range("F38").Value = 721
$ws.Range("F39").Value = 777
$ws.Range("F43").Value = 188
$ws.Range("F45").Value = 51
$ws.Range("F49").Value = 269
